# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the first-of-quarter date (e.g. 1988-07-01) used to key
# each quarterly GDP observation, but it should instead hold the
# "first-release" date for that quarter: the 15th of the month following
# the quarter-start month (e.g. 1988-07-01 -> 1988-08-15).
#
# Shift every date in column A, rows 2..150, forward by one calendar
# month and pin the day-of-month to 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldDate = [DateTime]::FromOADate($cell.Value2)
    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value2 = $newDate.ToOADate()
}
